$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.356900691986084
$ws.Range("B1").Value = 2.863547563552856
$ws.Range("C1").Value = 3.924045324325562
$ws.Range("D1").Value = 3.43274712562561
$ws.Range("E1").Value = 0.7882069945335388
